# Swap the data content between row 4 and row 5 (columns D, J, K, L, M, N, O, P, Q)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

foreach ($col in $cols) {
    $cell4 = $ws.Range($col + "4")
    $cell5 = $ws.Range($col + "5")
    $v4 = $cell4.Value2
    $v5 = $cell5.Value2
    $cell4.Value2 = $v5
    $cell5.Value2 = $v4
}
